$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old row 3 (the 3/26/2012 test), since it was a repeat of the
# 3/27/2012 testing and is being removed per the commit message.
[void]$ws.Rows(3).Delete()

# --- Row 2 (3/13/2012 test): "Buck Supply" -> "buck filtering" ---
$ws.Range("A2").Value = 40981
[void]$ws.Range("C2").Copy()
[void]$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("B2").Value = "buck filtering"
$ws.Range("D2").Value = "The capacitors were put in a low power buck supply to test for degredation"
$ws.Range("E2").Value = "0x0002"
$ws.Range("F2").Value = "0x0003"

# --- Row 3 (now 3/27/2012 test, formerly row 4): "Baseline measurement" -> "leakage" ---
$ws.Range("B3").Value = "leakage"
$ws.Range("D3").Value = "Measure leakage of tant and ti capacitors"
$ws.Range("E3").Value = "0x0003"
$ws.Range("F3").Value = "Tant - TAP475K006SCS"

# --- Row 4 (4/12/2012 test, formerly row 5): unchanged content, just shifted up ---
# (content already matches target after the row delete/shift, nothing to change)

$excel.CutCopyMode = 0

# Columns F and H now contain the wider "Tant - TAP475K006SCS" string, so
# widen them to fit (column G is untouched and keeps its original width).
$ws.Columns("F").ColumnWidth = 18.833333333333332
$ws.Columns("H").ColumnWidth = 18.833333333333332

# Update the selection shown in the worksheet view to match the final state.
[void]$ws.Range("H8:H9").Select()
